$wb = $excel.ActiveWorkbook

# Report that the handoff transform failed for the markdown source file, on
# both the zh-cn and de-de sheets: the status changes, the "latest handoff
# file" link/value is cleared (no file was produced), the "latest handoff
# datetime" reverts to the zero-value placeholder, and the handoff reason
# flips from "Include" to "Ignored".
#
# The two per-language sheets each keep two hyperlinks afterwards: the
# source markdown file (A2) and the .localization-config file (A3). The
# hyperlink that used to sit on the now-cleared C2 (the generated .xlf
# handoff file) goes away. This engine's `Range.Hyperlinks.Delete()` drops
# every hyperlink on the sheet rather than just the ones in the range, so
# the surviving two are recreated afterwards from their original
# target/display text.
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/0bf46882c663a8382a7dde475222c444aedc07dc/e2e/2546b087-f924-469b-bc49-93289baa5b90.md"
$mdDisplay = "2546b087-f924-469b-bc49-93289baa5b90.md"
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/0bf46882c663a8382a7dde475222c444aedc07dc/.localization-config"
$configDisplay = ".localization-config"

$sheetNames = @("zh-cn", "de-de")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Status column (B2): "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff Datetime (D2): reset to the "never happened" placeholder.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (H2): "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"

    # Drop every hyperlink on the sheet, then rebuild only the two that
    # should survive (A2, A3) -- this also takes care of removing the
    # hyperlink that was anchored on C2.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdTarget, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configTarget, "", "", $configDisplay)

    # Re-adding hyperlinks resets their cell style to the engine's default
    # hyperlink look; restore the workbook's original custom hyperlink font
    # (underlined cornflower blue Calibri 11) on those two cells.
    $hyperlinkCells = $ws.Range("A2:A3")
    $hyperlinkCells.Font.Name = "Calibri"
    $hyperlinkCells.Font.Size = 11
    $hyperlinkCells.Font.Underline = 2
    $hyperlinkCells.Font.Color = 15570276

    # Latest Handoff File (C2): remove the value entirely -- no file was
    # produced since the transform failed.
    $ws.Range("C2").Clear()
}

# The Overview sheet's per-language status cells share the same "status"
# string, so they follow the status text change too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"
